$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B (rows 2-7 existing tickers replaced, rows 8-17 newly filled) ---
$bValues = @(
    "NSE:20MICRONS",   # row 2
    "NSE:ANMOL",       # row 3
    "NSE:AVONMORE",    # row 4
    "NSE:BALMLAWRIE",  # row 5
    "NSE:CELLO",       # row 6
    "NSE:CRAFTSMAN",   # row 7
    "NSE:EMMBI",       # row 8
    "NSE:GAEL",        # row 9
    "NSE:GNFC",        # row 10
    "NSE:HARIOMPIPE",  # row 11
    "NSE:IPL",         # row 12
    "NSE:KDDL",        # row 13
    "NSE:LINDEINDIA",  # row 14
    "NSE:MMTC",        # row 15
    "NSE:MSTCLTD",     # row 16
    "NSE:NIPPOBATRY"   # row 17
)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# --- Update column C (rows 2-7 replaced with new tickers, rows 8-17 cleared) ---
$cValues = @(
    "NSE:DBL",        # row 2
    "NSE:ESAFSFB",    # row 3
    "NSE:HMAAGRO",    # row 4
    "NSE:JUSTDIAL",   # row 5
    "NSE:NPBET",      # row 6
    "NSE:ORIENTHOT"   # row 7
)
for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}
for ($row = 8; $row -le 17; $row++) {
    $ws.Cells.Item($row, 3).ClearContents()
}

# --- Update column E (row 2 ticker replaced) ---
$ws.Cells.Item(2, 5).Value = "NSE:COLPAL"

# --- Append new rows 18 and 19 ---
# Copy formatting from the last existing data row (17) for column A so the
# numbering style (bold, bordered, centered) is preserved.
$ws.Range("A17").Copy()
$ws.Range("A18:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(19, 1).Value = 17

$ws.Cells.Item(18, 2).Value = "NSE:OBEROIRLTY"
$ws.Cells.Item(19, 2).Value = "NSE:RADAAN"
